$wb = $excel.ActiveWorkbook

# --- Rename sheets (by position, matches workbook.xml <sheets> order) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509960905965006"
$wb.Worksheets.Item(2).Name = "NB_TO-1650996092133394"
$wb.Worksheets.Item(3).Name = "RS_TO-1650996092133394"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509960921813798"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509960922454135"

# --- Sheet 1 (GNG_TO) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509960905564651.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960905725052.csv"
$ws1.Range("B4").Value = "go_stims-16509960905725052.csv"
$ws1.Range("B5").Value = "GNG_stims-16509960905965006.csv"

# --- Sheet 2 (NB_TO) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16509960912373827.csv"
$ws2.Range("B3").Value = "ZB-match_0-16509960910773842.csv"
$ws2.Range("B4").Value = "OB-1650996091405415.csv"
$ws2.Range("B5").Value = "ZB-match_0-16509960907245052.csv"
$ws2.Range("B6").Value = "TB-16509960920693789.csv"
$ws2.Range("B7").Value = "ZB-match_0-16509960908453841.csv"
$ws2.Range("B8").Value = "TB-16509960917573824.csv"
$ws2.Range("B9").Value = "TB-16509960921174114.csv"
$ws2.Range("B10").Value = "OB-16509960912934108.csv"

# --- Sheet 3 (RS_TO) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (TOL_TO) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509960921494126.csv"
$ws4.Range("B3").Value = "ZM_stims-1650996092133394.csv"
$ws4.Range("B4").Value = "MM_stims-16509960921654134.csv"
$ws4.Range("B5").Value = "ZM_stims-16509960921494126.csv"
$ws4.Range("B6").Value = "MM_stims-16509960921813798.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960921654134.csv"

# --- Sheet 5 (vSAT_TO) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650996092229382.csv"
$ws5.Range("B3").Value = "SAT_stims-16509960921813798.csv"
$ws5.Range("B4").Value = "SAT_stims-16509960921974132.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650996092213418.csv"
